$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# 1) Two existing header cells get their week/date label text corrected.
# ---------------------------------------------------------------------------
$ws.Range("A20").Value = "日期：2018.10.10 第六周周三"
$ws.Range("A41").Value = "日期：2018.10.15 第七周周一"

# ---------------------------------------------------------------------------
# 2) New weekly block, rows 51-59, mirrors the rows 41-49 block directly
#    above it (same team roster / same layout), with new plan text for the
#    10/17 (week 7, Thursday) status update.
#
#    Merged ranges are merged FIRST and then formatted as a whole (rather
#    than formatting the individual cells and merging afterwards) so that
#    every covered cell keeps the sheet's single "thin box on all sides"
#    border style instead of Excel's neighbour-aware border splitting.
# ---------------------------------------------------------------------------

# -- Row 51: section/date header (bold, centered, merged A:D) --------------
$r51 = $ws.Range("A51:D51")
$r51.Merge()
$r51.Borders.LineStyle = 1
$r51.Font.Bold = $true
$r51.Font.Size = 10
$r51.HorizontalAlignment = -4108   # xlCenter
$r51.VerticalAlignment = -4108     # xlCenter
$ws.Range("A51").Value = "日期：2018.10.17 第七周周四"

# -- Row 52: column headers (bold) ------------------------------------------
$r52 = $ws.Range("A52:D52")
$r52.Borders.LineStyle = 1
$r52.Font.Bold = $true
$r52.Font.Size = 10
$r52.VerticalAlignment = -4108
$ws.Range("A52").Value = "组员"
$ws.Range("B52").Value = "计划内容"
$ws.Range("C52").Value = "完成情况"
$ws.Range("D52").Value = "备注"

# -- Rows 53-57: one row per team member ------------------------------------
$members  = @("王伟锋", "陈升云", "林玮成", "吴帅辰", "李海洋")
$plans    = @(
    "完成创建数据库表，完成注册servlet",
    "完成ui界面代码的编写（登录，注册）",
    "完成用例图的修改",
    "着手编写系统管理员代码",
    "编写网络交互文档"
)

for ($i = 0; $i -lt 5; $i++) {
    $row = 53 + $i
    $rowRange = $ws.Range("A$($row):D$($row)")
    $rowRange.Borders.LineStyle = 1
    $rowRange.VerticalAlignment = -4108

    $ws.Range("A$($row)").Value = $members[$i]

    $planCell = $ws.Range("B$($row)")
    $planCell.Value = $plans[$i]

    $statusCell = $ws.Range("C$($row)")
    $statusCell.Value = 1
    $statusCell.NumberFormat = "0%"

    $ws.Range("D$($row)").Value = $null
}

# -- Row 58-59: summary header, merged A58:D59 -------------------------------
$r58 = $ws.Range("A58:D59")
$r58.Merge()
$r58.Borders.LineStyle = 1
$r58.HorizontalAlignment = -4131   # xlLeft
$r58.VerticalAlignment = -4108
$ws.Range("A58").Value = "总结：熟悉了er图的设计流程,让数据库表创建变得更加简单"

# ---------------------------------------------------------------------------
# 3) Update the view/selection state to match where the editor ended up.
# ---------------------------------------------------------------------------
$ws.Range("C55").Select()
